# Update cryptos list values per latest data refresh (2024-02-23)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.141.41'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").Value = '2.946.76'
$ws.Range("E3").Value = '  -0.92%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'381.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("D6").Value = "'102.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.29%  '
$ws.Range("D7").Value = "'0.538"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.36%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = "'0.588"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.22%  '
$ws.Range("D10").Value = "'36.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.29%  '
$ws.Range("E11").Value = '  -0.72%  '
$ws.Range("D12").Value = "'0.0841"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("D13").Value = '3.419.53'
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("D14").Value = "'18.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.27%  '
$ws.Range("D15").Value = "'7.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.15%  '
$ws.Range("D16").Value = '2.950.30'
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("E17").Value = '  +2.94%  '
$ws.Range("D18").Value = '51.100.52'
$ws.Range("E18").Value = '  -1.42%  '
$ws.Range("D19").Value = "'3.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.61%  '
$ws.Range("D20").Value = "'7.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.08%  '
$ws.Range("D21").Value = "'12.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.80%  '
$ws.Range("D22").Value = '0.0₃0953'
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").Value = "'68.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D24").Value = "'262.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("D25").Value = "'2.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.57%  '
$ws.Range("D26").Value = "'8.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.96%  '
$ws.Range("D27").Value = "'7.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.61%  '
$ws.Range("D28").Value = "'0.169"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").Value = "'0.113"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.44%  '
$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").Value = "'25.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.53%  '
$ws.Range("D32").Value = "'9.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.20%  '
$ws.Range("D33").Value = "'34.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.89%  '
$ws.Range("E34").Value = '  +4.10%  '
$ws.Range("D35").Value = "'50.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.72%  '
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = "'2.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.99%  '
$ws.Range("D39").Value = "'16.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.48%  '
$ws.Range("D40").Value = "'2.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.45%  '
$ws.Range("D41").Value = "'0.115"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.46%  '
$ws.Range("D42").Value = "'1.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.99%  '
$ws.Range("D43").Value = "'120.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.30%  '
$ws.Range("D44").Value = "'21.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.85%  '
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("E46").Value = '  +2.89%  '
$ws.Range("E47").Value = '  -3.66%  '
$ws.Range("D48").Value = "'3.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("D49").Value = '2.009.83'
$ws.Range("E49").Value = '  -1.33%  '
$ws.Range("E50").Value = '  +5.06%  '
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").Value = "'0.480"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +14.60%  '
